$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("schedule")
$ws.Range("B2").Value = "Hello"
